$wb = $excel.ActiveWorkbook

# Update the variable load value on the "buses" sheet.
# B5 holds the literal value; B6:B13 reference it via formulas and will
# recompute automatically once B5 changes.
$wsBuses = $wb.Worksheets.Item("buses")
$wsBuses.Range("B5").Value = -0.1

# Make "buses" the active/selected sheet (previously "lines" was active).
$wsBuses.Activate()
$wsBuses.Select()

$wb.Save()
